# Update the crypto price/volume snapshot to the values scraped in the
# latest GitHub Actions run (Wed May 31 15:35:30 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.975.82"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "1.861.07"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'305.91"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.5069"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "'0.3736"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "'0.07122"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'0.8872"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "'20.51"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'0.07551"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "1.852.84"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'89.06"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'0.000008360"
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").Value = "'14.07"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.023.94"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "'5.054"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "2.092.22"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").Value = "'10.48"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  -1.69%  "
# Row 25: Monero/Toncoin swapped position in the ranking
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.846"
$ws.Range("E25").Value = "  +0.31%  "

# Row 26: Monero/Toncoin swapped position in the ranking
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'147.43"
$ws.Range("E26").Value = "  -3.77%  "

$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "'2.085"
$ws.Range("E28").Value = "  -4.45%  "
$ws.Range("D29").Value = "'112.61"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "'4.664"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "'4.641"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Value = "'0.09032"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "'0.05112"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").Value = "'3.057"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "'1.152"
$ws.Range("E35").Value = "  -5.56%  "
$ws.Range("D36").Value = "'0.7274"
$ws.Range("E36").Value = "  -6.47%  "
$ws.Range("D37").Value = "'0.02040"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").Value = "'3.039"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").Value = "'2.459"
$ws.Range("E39").Value = "  -5.63%  "
$ws.Range("D40").Value = "'1.071"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "'0.5316"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("D42").Value = "'6.575"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'115.56"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "'8.299"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").Value = "'0.1469"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'0.4607"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D48").Value = "'10.04"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").Value = "'1.562"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").Value = "'36.49"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'63.94"
$ws.Range("E51").Value = "  -3.92%  "
